$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert two new rows before the current row 7 ("admin" row), pushing it down to row 9.
$ws.Rows.Item(7).Insert()
$ws.Rows.Item(7).Insert()

# 2) The newly inserted rows 7 and 8 inherited row 6's formatting (border-less, font-3 text /
#    font-4 number style). Give them the same row heights as the other data-template row (19.5).
$ws.Rows.Item(7).RowHeight = 19.5
$ws.Rows.Item(8).RowHeight = 19.5

# 3) Swap the font colors that were mixed up between the two "Calibri, 11pt" fonts used by
#    row 6 (separator row) and row 7/8/9 (data rows): row 6's text should end up explicit black,
#    while the data rows should use the normal theme color.
$ws.Range("A6:C6").Font.Color = 0
$ws.Range("E6:F6").Font.Color = 0

for ($r = 7; $r -le 8; $r++) {
    $ws.Cells.Item($r, 1).Font.ThemeColor = 1
    $ws.Cells.Item($r, 2).Font.ThemeColor = 1
    $ws.Cells.Item($r, 3).Font.ThemeColor = 1
    $ws.Cells.Item($r, 5).Font.ThemeColor = 1
    $ws.Cells.Item($r, 6).Font.ThemeColor = 1
    $ws.Cells.Item($r, 4).Font.ThemeColor = 1
}

# 4) Rows 5 and 6 become fully blank placeholder rows (keep their style, drop their text).
$ws.Range("A5:F5").Value2 = ""
$ws.Range("A6:F6").Value2 = ""

# 5) Rows 7 and 8 are new blank "registration template" rows -- same shape as the data rows,
#    left empty for now.
$ws.Range("A7:C7").Value2 = ""
$ws.Range("E7:F7").Value2 = ""
$ws.Range("A8:C8").Value2 = ""
$ws.Range("E8:F8").Value2 = ""

# 6) Update the admin row (now row 9): the counter column moves from 0 to 13.
$ws.Cells.Item(9, 4).Value2 = 13

Write-Host "Done"
